$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = 42605.886689814812
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B4").Value = -36
$ws.Range("C4").Value = 43
$ws.Range("D4").Value = 56
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 22079
$ws.Range("H4").Value = 5102
$ws.Range("I4").Value = 323
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 52
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = "Named"
